$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 17.956883
$ws.Range("H2").Value = 35.913766
$ws.Range("I2").Value = 0.3392380274206944
$ws.Range("J2").Value = 0.2584869083704147
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.046576
$ws.Range("N2").Value = 0.093152
$ws.Range("O2").Value = 0.07994370161857847
$ws.Range("P2").Value = 0.07994370161857847
$ws.Range("Q2").Value = 0.8363597826079999
$ws.Range("R2").Value = 3.345439130432
$ws.Range("S2").Value = 0.02711994364179514
$ws.Range("T2").Value = 0.02066440027507326
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 17.956883
$ws.Range("H3").Value = 35.913766
$ws.Range("I3").Value = 0.3392380274206944
$ws.Range("J3").Value = 0.2584869083704147
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.536034
$ws.Range("N3").Value = 1.072068
$ws.Range("O3").Value = 0.9200562983814217
$ws.Range("P3").Value = 0.9200562983814217
$ws.Range("Q3").Value = 9.625499822021998
$ws.Range("R3").Value = 38.50199928808799
$ws.Range("S3").Value = 0.3121180837788993
$ws.Range("T3").Value = 0.2378225080953414
$ws.Range("D4").Value = "ECs"
$ws.Range("I4").Value = 0.6187742881378531
$ws.Range("J4").Value = 0.7072248972319991
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.046576
$ws.Range("N4").Value = 0.093152
$ws.Range("O4").Value = 0.07994370161857847
$ws.Range("P4").Value = 0.07994370161857847
$ws.Range("Q4").Value = 1.525530416048
$ws.Range("R4").Value = 9.153182496288
$ws.Range("S4").Value = 0.04946710706014083
$ws.Range("T4").Value = 0.05653817616154476
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.6187742881378531
$ws.Range("J5").Value = 0.7072248972319991
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.536034
$ws.Range("N5").Value = 1.072068
$ws.Range("O5").Value = 0.9200562983814217
$ws.Range("P5").Value = 0.9200562983814217
$ws.Range("Q5").Value = 17.557028749482
$ws.Range("R5").Value = 105.342172496892
$ws.Range("S5").Value = 0.5693071810777123
$ws.Range("T5").Value = 0.6506867210704544
$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 1.903653
$ws.Range("H6").Value = 3.807306
$ws.Range("I6").Value = 0.03596345137480081
$ws.Range("J6").Value = 0.02740282812891664
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.046576
$ws.Range("N6").Value = 0.093152
$ws.Range("O6").Value = 0.07994370161857847
$ws.Range("P6").Value = 0.07994370161857847
$ws.Range("Q6").Value = 0.088664542128
$ws.Range("R6").Value = 0.354658168512
$ws.Range("S6").Value = 0.002875051425881331
$ws.Range("T6").Value = 0.002190683515443301
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 1.903653
$ws.Range("H7").Value = 3.807306
$ws.Range("I7").Value = 0.03596345137480081
$ws.Range("J7").Value = 0.02740282812891664
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.536034
$ws.Range("N7").Value = 1.072068
$ws.Range("O7").Value = 0.9200562983814217
$ws.Range("P7").Value = 0.9200562983814217
$ws.Range("Q7").Value = 1.020422732202
$ws.Range("R7").Value = 4.081690928808
$ws.Range("S7").Value = 0.03308839994891948
$ws.Range("T7").Value = 0.02521214461347334
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3188806666666666
$ws.Range("H8").Value = 0.956642
$ws.Range("I8").Value = 0.006024233066651711
$ws.Range("J8").Value = 0.006885366268669519
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.046576
$ws.Range("N8").Value = 0.093152
$ws.Range("O8").Value = 0.07994370161857847
$ws.Range("P8").Value = 0.07994370161857847
$ws.Range("Q8").Value = 0.01485218593066667
$ws.Range("R8").Value = 0.08911311558399999
$ws.Range("S8").Value = 0.0004815994907611784
$ws.Range("T8").Value = 0.0005504416665171411
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3188806666666666
$ws.Range("H9").Value = 0.956642
$ws.Range("I9").Value = 0.006024233066651711
$ws.Range("J9").Value = 0.006885366268669519
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.536034
$ws.Range("N9").Value = 1.072068
$ws.Range("O9").Value = 0.9200562983814217
$ws.Range("P9").Value = 0.9200562983814217
$ws.Range("Q9").Value = 0.170930879276
$ws.Range("R9").Value = 1.025585275656
$ws.Range("S9").Value = 0.005542633575890534
$ws.Range("T9").Value = 0.00633492460215238
